$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 524.1111
$ws.Range("I33").Value = 411.77777
$ws.Range("K33").Value = 411.77777
$ws.Range("M33").Value = -182.77777
$ws.Range("H38").Value = 1447.2
$ws.Range("I38").Value = 176.81818
$ws.Range("J38").Value = 2999.889
$ws.Range("K38").Value = 530.4545400000001
$ws.Range("L38").Value = 8999.667000000001
$ws.Range("M38").Value = -158.4545400000001
$ws.Range("N38").Value = -9743.667000000001
$ws.Range("H58").Value = 23019.896
$ws.Range("J58").Value = 26258.096
$ws.Range("L58").Value = 78774.288
$ws.Range("N58").Value = -79074.288
$ws.Range("H96").Value = 748.2308
$ws.Range("I96").Value = 523.3333
$ws.Range("J96").Value = 941
$ws.Range("K96").Value = 1569.9999
$ws.Range("L96").Value = 2823
$ws.Range("M96").Value = -196.9999
$ws.Range("N96").Value = -5569
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 1312.71
$ws.Range("I112").Value = 547.75
$ws.Range("J112").Value = 1344.5834
$ws.Range("K112").Value = 1643.25
$ws.Range("L112").Value = 4033.7502
$ws.Range("M112").Value = -535.25
$ws.Range("N112").Value = -6249.7502
$ws.Range("H132").Value = 3705922.8
$ws.Range("I132").Value = 4168612.8
$ws.Range("J132").Value = 4401.8335
$ws.Range("K132").Value = 12505838.4
$ws.Range("L132").Value = 13205.5005
$ws.Range("M132").Value = -12503308.4
$ws.Range("N132").Value = -18265.5005
$ws.Range("H137").Value = 2448.7058
$ws.Range("I137").Value = 2283.628
$ws.Range("K137").Value = 6850.884
$ws.Range("M137").Value = -4300.884
$ws.Range("H138").Value = 1874.93
$ws.Range("I138").Value = 774.175
$ws.Range("J138").Value = 2608.7666
$ws.Range("K138").Value = 2322.525
$ws.Range("L138").Value = 7826.2998
$ws.Range("M138").Value = 2817.475
$ws.Range("N138").Value = -18106.2998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1981.46
$ws.Range("I61").Value = 793.8
$ws.Range("J61").Value = 4752.6665
$ws.Range("K61").Value = 793.8
$ws.Range("L61").Value = 4752.6665
$ws.Range("M61").Value = -581.8
$ws.Range("N61").Value = -5176.6665
$ws.Range("H74").Value = 784.62067
$ws.Range("I74").Value = 692.7273
$ws.Range("J74").Value = 1073.4286
$ws.Range("K74").Value = 692.7273
$ws.Range("L74").Value = 1073.4286
$ws.Range("M74").Value = 181.2727
$ws.Range("N74").Value = -2821.4286
$ws.Range("H77").Value = 784.62067
$ws.Range("I77").Value = 692.7273
$ws.Range("J77").Value = 1073.4286
$ws.Range("K77").Value = 3463.6365
$ws.Range("L77").Value = 5367.143
$ws.Range("M77").Value = 904.3634999999999
$ws.Range("N77").Value = -14103.143
$ws.Range("H132").Value = 2110.9473
$ws.Range("I132").Value = 1454.8
$ws.Range("J132").Value = 3154.818
$ws.Range("K132").Value = 4364.4
$ws.Range("L132").Value = 9464.454000000002
$ws.Range("M132").Value = -1834.4
$ws.Range("N132").Value = -14524.454
$ws.Range("H136").Value = 1981.46
$ws.Range("I136").Value = 793.8
$ws.Range("J136").Value = 4752.6665
$ws.Range("K136").Value = 2381.4
$ws.Range("L136").Value = 14257.9995
$ws.Range("M136").Value = 168.6000000000004
$ws.Range("N136").Value = -19357.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 30779.523
$ws.Range("J133").Value = 30779.523
$ws.Range("L133").Value = 30779.523
$ws.Range("N133").Value = -40899.523
$ws.Range("H134").Value = 1797.6364
$ws.Range("I134").Value = 1449.3182
$ws.Range("J134").Value = 3190.9092
$ws.Range("K134").Value = 4347.9546
$ws.Range("L134").Value = 9572.7276
$ws.Range("M134").Value = -1812.9546
$ws.Range("N134").Value = -14642.7276

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3423.3684
$ws.Range("I31").Value = 2317.4443
$ws.Range("J31").Value = 4418.7
$ws.Range("K31").Value = 2317.4443
$ws.Range("L31").Value = 4418.7
$ws.Range("M31").Value = -2022.4443
$ws.Range("N31").Value = -5008.7
$ws.Range("H34").Value = 3423.3684
$ws.Range("I34").Value = 2317.4443
$ws.Range("J34").Value = 4418.7
$ws.Range("K34").Value = 2317.4443
$ws.Range("L34").Value = 4418.7
$ws.Range("M34").Value = -2115.4443
$ws.Range("N34").Value = -4822.7
$ws.Range("H58").Value = 7938776
$ws.Range("I58").Value = 1332.7805
$ws.Range("J58").Value = 22731284
$ws.Range("K58").Value = 1332.7805
$ws.Range("L58").Value = 22731284
$ws.Range("M58").Value = -1129.7805
$ws.Range("N58").Value = -22731690
$ws.Range("H127").Value = 32993.08
$ws.Range("J127").Value = 32993.08
$ws.Range("L127").Value = 32993.08
$ws.Range("N127").Value = -42913.08
$ws.Range("H132").Value = 1842.7727
$ws.Range("I132").Value = 1395.0667
$ws.Range("J132").Value = 2802.1428
$ws.Range("K132").Value = 4185.2001
$ws.Range("L132").Value = 8406.428400000001
$ws.Range("M132").Value = -1655.2001
$ws.Range("N132").Value = -13466.4284
$ws.Range("H134").Value = 1998.6428
$ws.Range("I134").Value = 894.64703
$ws.Range("J134").Value = 3704.818
$ws.Range("K134").Value = 2683.94109
$ws.Range("L134").Value = 11114.454
$ws.Range("M134").Value = -148.9410899999998
$ws.Range("N134").Value = -16184.454
$ws.Range("H136").Value = 7938776
$ws.Range("I136").Value = 1332.7805
$ws.Range("J136").Value = 22731284
$ws.Range("K136").Value = 3998.3415
$ws.Range("L136").Value = 68193852
$ws.Range("M136").Value = -1448.3415
$ws.Range("N136").Value = -68198952
$ws.Range("H137").Value = 29612.5
$ws.Range("J137").Value = 29612.5
$ws.Range("L137").Value = 29612.5
$ws.Range("N137").Value = -39812.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 933.3333
$ws.Range("J97").Value = 1088.5714
$ws.Range("L97").Value = 3265.7142
$ws.Range("N97").Value = -4257.7142
$ws.Range("H131").Value = 1442.1818
$ws.Range("I131").Value = 1887.1428
$ws.Range("K131").Value = 5661.428400000001
$ws.Range("M131").Value = -621.4284000000007
$ws.Range("H136").Value = 2189.818
$ws.Range("I136").Value = 1488.8235
$ws.Range("J136").Value = 4573.2
$ws.Range("K136").Value = 4466.470499999999
$ws.Range("L136").Value = 13719.6
$ws.Range("M136").Value = 633.5295000000006
$ws.Range("N136").Value = -23919.6
$ws.Range("H138").Value = 4664.1
$ws.Range("J138").Value = 6783.3335
$ws.Range("L138").Value = 20350.0005
$ws.Range("N138").Value = -30630.0005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 27000
$ws.Range("J140").Value = 27000
$ws.Range("L140").Value = 27000
$ws.Range("N140").Value = -37360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 274875.12
$ws.Range("J2").Value = 49750.25
$ws.Range("L2").Value = 49750.25
$ws.Range("N2").Value = -49974.25
$ws.Range("H22").Value = 76924410
$ws.Range("I22").Value = 200000530
$ws.Range("J22").Value = 1835.5
$ws.Range("K22").Value = 200000530
$ws.Range("L22").Value = 1835.5
$ws.Range("M22").Value = -200000235
$ws.Range("N22").Value = -2425.5
$ws.Range("H24").Value = 37503.5
$ws.Range("J24").Value = 37503.5
$ws.Range("L24").Value = 37503.5
$ws.Range("N24").Value = -38189.5
$ws.Range("H27").Value = 76924410
$ws.Range("I27").Value = 200000530
$ws.Range("J27").Value = 1835.5
$ws.Range("K27").Value = 200000530
$ws.Range("L27").Value = 1835.5
$ws.Range("M27").Value = -200000423
$ws.Range("N27").Value = -2049.5
$ws.Range("H132").Value = 2402.6316
$ws.Range("I132").Value = 1473.5652
$ws.Range("J132").Value = 3827.2
$ws.Range("K132").Value = 4420.6956
$ws.Range("L132").Value = 11481.6
$ws.Range("M132").Value = -1890.6956
$ws.Range("N132").Value = -16541.6
$ws.Range("H136").Value = 2023.5385
$ws.Range("I136").Value = 1201.2222
$ws.Range("J136").Value = 3873.75
$ws.Range("K136").Value = 3603.6666
$ws.Range("L136").Value = 11621.25
$ws.Range("M136").Value = -1053.6666
$ws.Range("N136").Value = -16721.25
$ws.Range("H137").Value = 29565.385
$ws.Range("J137").Value = 29565.385
$ws.Range("L137").Value = 29565.385
$ws.Range("N137").Value = -39765.38499999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43979.5
$ws.Range("J46").Value = 43979.5
$ws.Range("L46").Value = 43979.5
$ws.Range("N46").Value = -44441.5
$ws.Range("H122").Value = 478584.94
$ws.Range("I122").Value = 668657.4
$ws.Range("J122").Value = 3403.8333
$ws.Range("K122").Value = 2005972.2
$ws.Range("L122").Value = 10211.4999
$ws.Range("M122").Value = -2003522.2
$ws.Range("N122").Value = -15111.4999
$ws.Range("H132").Value = 12715.44
$ws.Range("I132").Value = 2859.543
$ws.Range("K132").Value = 8578.629000000001
$ws.Range("M132").Value = -6048.629000000001
$ws.Range("H134").Value = 43979.5
$ws.Range("J134").Value = 43979.5
$ws.Range("L134").Value = 131938.5
$ws.Range("N134").Value = -137008.5
$ws.Range("H136").Value = 1109.5581
$ws.Range("I136").Value = 775.4706
$ws.Range("J136").Value = 2371.6667
$ws.Range("K136").Value = 2326.4118
$ws.Range("L136").Value = 7115.000100000001
$ws.Range("M136").Value = 223.5882000000001
$ws.Range("N136").Value = -12215.0001
$ws.Range("H140").Value = 79999
$ws.Range("J140").Value = 79999
$ws.Range("L140").Value = 79999
$ws.Range("N140").Value = -90359
